# Notulen base.xlsx - "Steven map hernieuwd (alweer)"
# Resets the meeting-minutes template: clears the previously filled-in
# agenda/notes content, updates a couple of labels, and makes the agenda
# point column bold again.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Title cell: add the trailing "fill in the date" dots (ellipsis x4 + ..) ---
$ellipsis = [char]0x2026
$ws.Range("A1").Value = "Verslag overleg op " + "$ellipsis$ellipsis$ellipsis$ellipsis.."

# --- Agenda point labels ---
# Point 3 is renamed from a reference to last meeting's minutes to a
# generic "determine next meeting item" placeholder.
$ws.Range("A9").Value = "3 Vaststellen volgende vergaderpunt"

# --- Clear stale, meeting-specific content so the sheet is reusable ---
$ws.Range("B7").ClearContents()
$ws.Range("A14").ClearContents()
$ws.Range("A15").ClearContents()
$ws.Range("A16").ClearContents()
$ws.Range("A17:A22").ClearContents()

# --- Make the agenda-point column (A7:A22) bold again ---
$ws.Range("A7:A22").Font.Bold = $true

# --- Rows that no longer need their custom (wrapped-text) height fall
#     back to the sheet's default row height ---
$ws.Range("A10:A15").EntireRow.AutoFit()
$ws.Range("A21:A22").EntireRow.AutoFit()

# --- Selection / scroll position: the view now opens at the top of the
#     sheet (no frozen topLeftCell) with A11 selected ---
$ws.Range("A11").Select() | Out-Null
